# "TestData add 3 colums" -- add a 3-column header row (TestCase / UserName /
# Pass) to Sheet1 and widen column B to fit "UserName".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # Sheet1 is already the active/selected sheet

$ws.Range("A1").Value = "TestCase"
$ws.Range("B1").Value = "UserName"
$ws.Range("C1").Value = "Pass"

# Widen column B (~15.43 stored char-width in the source file) to fit the
# "UserName" header.
$ws.Columns.Item(2).ColumnWidth = 14.7115

# Leave the selection on C1, matching the saved view state.
$ws.Range("C1").Select()
